# Weekly fruit/vegetable price log: insert a newly-recorded sampling date
# (row 257) and shift the remaining rows down one position.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 257; Excel shifts rows 257..331 down to 258..332
# and carries the date column's number format (style id 2) onto the new row.
$ws.Rows.Item(257).Insert()

# Populate the newly inserted row 257 with the new observation.
$ws.Cells.Item(257, 1).Value = 3
$ws.Cells.Item(257, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(257, 3).Value = "Coquimbo"
$ws.Cells.Item(257, 4).Value = 44736
$ws.Cells.Item(257, 5).Value = 5
$ws.Cells.Item(257, 6).Value = 100112039
$ws.Cells.Item(257, 7).Value = "Ciboulette"
$ws.Cells.Item(257, 8).Value = "Sin especificar"
$ws.Cells.Item(257, 9).Value = "Primera"
$ws.Cells.Item(257, 10).Value = 120
$ws.Cells.Item(257, 11).Value = 1500
$ws.Cells.Item(257, 12).Value = 1500
$ws.Cells.Item(257, 13).Value = 1500
$ws.Cells.Item(257, 14).Value = "`$/docena de atados"
$ws.Cells.Item(257, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(257, 16).Value = 500
$ws.Cells.Item(257, 17).Value = 3
$ws.Cells.Item(257, 18).Value = "Hortaliza"
